# Auto-generated edit script: update profit-calculation values across all sheets
# per the scheduled-runner price refresh (Typhon_Profits workbook).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 8236  # H19: 10250 -> 8236
$ws.Cells.Item(19, 10).Value = 726.6667  # J19: 1000 -> 726.6667
$ws.Cells.Item(19, 12).Value = 726.6667  # L19: 1000 -> 726.6667
$ws.Cells.Item(19, 14).Value = -1076.6667  # N19: -1350 -> -1076.6667
$ws.Cells.Item(40, 8).Value = 2275  # H40: 2457.1428 -> 2275
$ws.Cells.Item(40, 10).Value = 2171.4285  # J40: 2366.6667 -> 2171.4285
$ws.Cells.Item(40, 12).Value = 2171.4285  # L40: 2366.6667 -> 2171.4285
$ws.Cells.Item(40, 14).Value = -2521.4285  # N40: -2716.6667 -> -2521.4285
$ws.Cells.Item(53, 8).Value = 1800.6666  # H53: 144.8421 -> 1800.6666
$ws.Cells.Item(53, 9).Value = 156  # I53: 170 -> 156
$ws.Cells.Item(53, 10).Value = 2433.2307  # J53: 138.13333 -> 2433.2307
$ws.Cells.Item(53, 11).Value = 156  # K53: 170 -> 156
$ws.Cells.Item(53, 12).Value = 2433.2307  # L53: 138.13333 -> 2433.2307
$ws.Cells.Item(53, 13).Value = 481  # M53: 467 -> 481
$ws.Cells.Item(53, 14).Value = -3707.2307  # N53: -1412.13333 -> -3707.2307
$ws.Cells.Item(132, 8).Value = 2692.8918  # H132: 2760.75 -> 2692.8918
$ws.Cells.Item(132, 9).Value = 2625.182  # I132: 2699.4062 -> 2625.182
$ws.Cells.Item(132, 11).Value = 7875.545999999999  # K132: 8098.2186 -> 7875.545999999999
$ws.Cells.Item(132, 13).Value = -5345.545999999999  # M132: -5568.2186 -> -5345.545999999999
$ws.Cells.Item(138, 8).Value = 1654.247  # H138: 1682.3636 -> 1654.247
$ws.Cells.Item(138, 9).Value = 936.6667  # I138: 937.5333000000001 -> 936.6667
$ws.Cells.Item(138, 10).Value = 2551.2222  # J138: 2729.7812 -> 2551.2222
$ws.Cells.Item(138, 11).Value = 2810.0001  # K138: 2812.5999 -> 2810.0001
$ws.Cells.Item(138, 12).Value = 7653.6666  # L138: 8189.3436 -> 7653.6666
$ws.Cells.Item(138, 13).Value = 2329.9999  # M138: 2327.4001 -> 2329.9999
$ws.Cells.Item(138, 14).Value = -17933.6666  # N138: -18469.3436 -> -17933.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1393.4667  # H2: 1263.8572 -> 1393.4667
$ws.Cells.Item(2, 9).Value = 1324.6  # I2: 1120.92 -> 1324.6
$ws.Cells.Item(2, 10).Value = 1531.2  # J2: 1621.2 -> 1531.2
$ws.Cells.Item(2, 11).Value = 1324.6  # K2: 1120.92 -> 1324.6
$ws.Cells.Item(2, 12).Value = 1531.2  # L2: 1621.2 -> 1531.2
$ws.Cells.Item(2, 13).Value = -1211.6  # M2: -1007.92 -> -1211.6
$ws.Cells.Item(2, 14).Value = -1757.2  # N2: -1847.2 -> -1757.2
$ws.Cells.Item(34, 8).Value = 17000  # H34: 2222 -> 17000
$ws.Cells.Item(34, 9).Value = 0  # I34: 2222 -> 0
$ws.Cells.Item(34, 10).Value = 17000  # J34: 0 -> 17000
$ws.Cells.Item(34, 11).Value = 0  # K34: 2222 -> 0
$ws.Cells.Item(34, 12).Value = 17000  # L34: 0 -> 17000
$ws.Cells.Item(34, 13).ClearContents()  # M34: remove (was -1951)
$ws.Cells.Item(34, 14).Value = -17542  # N34: None -> -17542
$ws.Cells.Item(45, 8).Value = 2668.7222  # H45: 3014.3333 -> 2668.7222
$ws.Cells.Item(45, 9).Value = 3415.5  # I45: 4250.1665 -> 3415.5
$ws.Cells.Item(45, 10).Value = 2071.3  # J45: 2190.4443 -> 2071.3
$ws.Cells.Item(45, 11).Value = 3415.5  # K45: 4250.1665 -> 3415.5
$ws.Cells.Item(45, 12).Value = 2071.3  # L45: 2190.4443 -> 2071.3
$ws.Cells.Item(45, 13).Value = -3038.5  # M45: -3873.1665 -> -3038.5
$ws.Cells.Item(45, 14).Value = -2825.3  # N45: -2944.4443 -> -2825.3
$ws.Cells.Item(61, 8).Value = 2965.6365  # H61: 2504.0344 -> 2965.6365
$ws.Cells.Item(61, 9).Value = 2240.25  # I61: 2095.2222 -> 2240.25
$ws.Cells.Item(61, 10).Value = 4900  # J61: 3173 -> 4900
$ws.Cells.Item(61, 11).Value = 2240.25  # K61: 2095.2222 -> 2240.25
$ws.Cells.Item(61, 12).Value = 4900  # L61: 3173 -> 4900
$ws.Cells.Item(61, 13).Value = -2028.25  # M61: -1883.2222 -> -2028.25
$ws.Cells.Item(61, 14).Value = -5324  # N61: -3597 -> -5324
$ws.Cells.Item(74, 8).Value = 166667740  # H74: 55556370 -> 166667740
$ws.Cells.Item(74, 9).Value = 500000400  # I74: 111111530 -> 500000400
$ws.Cells.Item(74, 10).Value = 1407  # J74: 1209.2222 -> 1407
$ws.Cells.Item(74, 11).Value = 500000400  # K74: 111111530 -> 500000400
$ws.Cells.Item(74, 12).Value = 1407  # L74: 1209.2222 -> 1407
$ws.Cells.Item(74, 13).Value = -499999526  # M74: -111110656 -> -499999526
$ws.Cells.Item(74, 14).Value = -3155  # N74: -2957.2222 -> -3155
$ws.Cells.Item(77, 8).Value = 166667740  # H77: 55556370 -> 166667740
$ws.Cells.Item(77, 9).Value = 500000400  # I77: 111111530 -> 500000400
$ws.Cells.Item(77, 10).Value = 1407  # J77: 1209.2222 -> 1407
$ws.Cells.Item(77, 11).Value = 2500002000  # K77: 555557650 -> 2500002000
$ws.Cells.Item(77, 12).Value = 7035  # L77: 6046.111 -> 7035
$ws.Cells.Item(77, 13).Value = -2499997632  # M77: -555553282 -> -2499997632
$ws.Cells.Item(77, 14).Value = -15771  # N77: -14782.111 -> -15771
$ws.Cells.Item(116, 8).Value = 1393.4667  # H116: 1263.8572 -> 1393.4667
$ws.Cells.Item(116, 9).Value = 1324.6  # I116: 1120.92 -> 1324.6
$ws.Cells.Item(116, 10).Value = 1531.2  # J116: 1621.2 -> 1531.2
$ws.Cells.Item(116, 11).Value = 1324.6  # K116: 1120.92 -> 1324.6
$ws.Cells.Item(116, 12).Value = 1531.2  # L116: 1621.2 -> 1531.2
$ws.Cells.Item(116, 13).Value = 969.4000000000001  # M116: 1173.08 -> 969.4000000000001
$ws.Cells.Item(116, 14).Value = -6119.2  # N116: -6209.2 -> -6119.2
$ws.Cells.Item(136, 8).Value = 2965.6365  # H136: 2504.0344 -> 2965.6365
$ws.Cells.Item(136, 9).Value = 2240.25  # I136: 2095.2222 -> 2240.25
$ws.Cells.Item(136, 10).Value = 4900  # J136: 3173 -> 4900
$ws.Cells.Item(136, 11).Value = 6720.75  # K136: 6285.6666 -> 6720.75
$ws.Cells.Item(136, 12).Value = 14700  # L136: 9519 -> 14700
$ws.Cells.Item(136, 13).Value = -4170.75  # M136: -3735.6666 -> -4170.75
$ws.Cells.Item(136, 14).Value = -19800  # N136: -14619 -> -19800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1393.4667  # H3: 1263.8572 -> 1393.4667
$ws.Cells.Item(3, 9).Value = 1324.6  # I3: 1120.92 -> 1324.6
$ws.Cells.Item(3, 10).Value = 1531.2  # J3: 1621.2 -> 1531.2
$ws.Cells.Item(3, 11).Value = 1324.6  # K3: 1120.92 -> 1324.6
$ws.Cells.Item(3, 12).Value = 1531.2  # L3: 1621.2 -> 1531.2
$ws.Cells.Item(3, 13).Value = -1210.6  # M3: -1006.92 -> -1210.6
$ws.Cells.Item(3, 14).Value = -1759.2  # N3: -1849.2 -> -1759.2
$ws.Cells.Item(134, 8).Value = 73975.734  # H134: 43092.81 -> 73975.734
$ws.Cells.Item(134, 9).Value = 91761.336  # I134: 48344.043 -> 91761.336
$ws.Cells.Item(134, 11).Value = 275284.008  # K134: 145032.129 -> 275284.008
$ws.Cells.Item(134, 13).Value = -272749.008  # M134: -142497.129 -> -272749.008

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 50.285713  # H7: 52.615383 -> 50.285713
$ws.Cells.Item(7, 9).Value = 59.444443  # I7: 64.375 -> 59.444443
$ws.Cells.Item(7, 11).Value = 59.444443  # K7: 64.375 -> 59.444443
$ws.Cells.Item(7, 13).Value = 53.555557  # M7: 48.625 -> 53.555557
$ws.Cells.Item(31, 8).Value = 9843.41  # H31: 9741.737999999999 -> 9843.41
$ws.Cells.Item(31, 9).Value = 12567.038  # I31: 12138.034 -> 12567.038
$ws.Cells.Item(31, 11).Value = 12567.038  # K31: 12138.034 -> 12567.038
$ws.Cells.Item(31, 13).Value = -12272.038  # M31: -11843.034 -> -12272.038
$ws.Cells.Item(34, 8).Value = 9843.41  # H34: 9741.737999999999 -> 9843.41
$ws.Cells.Item(34, 9).Value = 12567.038  # I34: 12138.034 -> 12567.038
$ws.Cells.Item(34, 11).Value = 12567.038  # K34: 12138.034 -> 12567.038
$ws.Cells.Item(34, 13).Value = -12365.038  # M34: -11936.034 -> -12365.038
$ws.Cells.Item(60, 8).Value = 11442.263  # H60: 11911.333 -> 11442.263
$ws.Cells.Item(60, 9).Value = 3999.5  # I60: 5000 -> 3999.5
$ws.Cells.Item(60, 11).Value = 3999.5  # K60: 5000 -> 3999.5
$ws.Cells.Item(60, 13).Value = -3488.5  # M60: -4489 -> -3488.5
$ws.Cells.Item(134, 8).Value = 1220.6086  # H134: 1193.6459 -> 1220.6086
$ws.Cells.Item(134, 9).Value = 927.56525  # I134: 915.2917 -> 927.56525
$ws.Cells.Item(134, 10).Value = 1513.6522  # J134: 1472 -> 1513.6522
$ws.Cells.Item(134, 11).Value = 2782.69575  # K134: 2745.8751 -> 2782.69575
$ws.Cells.Item(134, 12).Value = 4540.9566  # L134: 4416 -> 4540.9566
$ws.Cells.Item(134, 13).Value = -247.6957499999999  # M134: -210.8751000000002 -> -247.6957499999999
$ws.Cells.Item(134, 14).Value = -9610.9566  # N134: -9486 -> -9610.9566

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 94  # H12: 101.85714 -> 94
$ws.Cells.Item(12, 9).Value = 35.5  # I12: 46.666668 -> 35.5
$ws.Cells.Item(12, 10).Value = 123.25  # J12: 116.90909 -> 123.25
$ws.Cells.Item(12, 11).Value = 106.5  # K12: 140.000004 -> 106.5
$ws.Cells.Item(12, 12).Value = 369.75  # L12: 350.72727 -> 369.75
$ws.Cells.Item(12, 13).Value = 66.5  # M12: 32.99999600000001 -> 66.5
$ws.Cells.Item(12, 14).Value = -715.75  # N12: -696.7272700000001 -> -715.75
$ws.Cells.Item(38, 8).Value = 62500116  # H38: 82.40000000000001 -> 62500116
$ws.Cells.Item(38, 9).Value = 86  # I38: 82.40000000000001 -> 86
$ws.Cells.Item(38, 10).Value = 166666830  # J38: 0 -> 166666830
$ws.Cells.Item(38, 11).Value = 258  # K38: 247.2 -> 258
$ws.Cells.Item(38, 12).Value = 500000490  # L38: 0 -> 500000490
$ws.Cells.Item(38, 13).Value = 89  # M38: 99.79999999999998 -> 89
$ws.Cells.Item(38, 14).Value = -500001184  # N38: None -> -500001184
$ws.Cells.Item(62, 8).Value = 7320.5454  # H62: 6966.4443 -> 7320.5454
$ws.Cells.Item(62, 10).Value = 9204.75  # J62: 9301.666999999999 -> 9204.75
$ws.Cells.Item(62, 12).Value = 27614.25  # L62: 27905.001 -> 27614.25
$ws.Cells.Item(62, 14).Value = -28986.25  # N62: -29277.001 -> -28986.25
$ws.Cells.Item(65, 8).Value = 7320.5454  # H65: 6966.4443 -> 7320.5454
$ws.Cells.Item(65, 10).Value = 9204.75  # J65: 9301.666999999999 -> 9204.75
$ws.Cells.Item(65, 12).Value = 82842.75  # L65: 83715.003 -> 82842.75
$ws.Cells.Item(65, 14).Value = -89706.75  # N65: -90579.003 -> -89706.75
$ws.Cells.Item(123, 8).Value = 4477.5  # H123: 4502.5 -> 4477.5
$ws.Cells.Item(123, 10).Value = 7445  # J123: 7495 -> 7445
$ws.Cells.Item(123, 12).Value = 22335  # L123: 22485 -> 22335
$ws.Cells.Item(123, 14).Value = -27235  # N123: -27385 -> -27235
$ws.Cells.Item(130, 8).Value = 400  # H130: 475 -> 400
$ws.Cells.Item(130, 9).Value = 400  # I130: 475 -> 400
$ws.Cells.Item(130, 11).Value = 1200  # K130: 1425 -> 1200
$ws.Cells.Item(130, 13).Value = 3820  # M130: 3595 -> 3820
$ws.Cells.Item(131, 8).Value = 734.08  # H131: 738.59 -> 734.08
$ws.Cells.Item(131, 10).Value = 746.90424  # J131: 751.70215 -> 746.90424
$ws.Cells.Item(131, 12).Value = 2240.71272  # L131: 2255.10645 -> 2240.71272
$ws.Cells.Item(131, 14).Value = -12320.71272  # N131: -12335.10645 -> -12320.71272

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2089.2856  # H22: 2089.3572 -> 2089.2856
$ws.Cells.Item(22, 9).Value = 1575  # I22: 1575.1 -> 1575
$ws.Cells.Item(22, 11).Value = 1575  # K22: 1575.1 -> 1575
$ws.Cells.Item(22, 13).Value = -1280  # M22: -1280.1 -> -1280
$ws.Cells.Item(27, 8).Value = 2089.2856  # H27: 2089.3572 -> 2089.2856
$ws.Cells.Item(27, 9).Value = 1575  # I27: 1575.1 -> 1575
$ws.Cells.Item(27, 11).Value = 1575  # K27: 1575.1 -> 1575
$ws.Cells.Item(27, 13).Value = -1468  # M27: -1468.1 -> -1468
$ws.Cells.Item(31, 8).Value = 932.4286  # H31: 812 -> 932.4286
$ws.Cells.Item(31, 9).Value = 702  # I31: 765 -> 702
$ws.Cells.Item(31, 10).Value = 1508.5  # J31: 1000 -> 1508.5
$ws.Cells.Item(31, 11).Value = 702  # K31: 765 -> 702
$ws.Cells.Item(31, 12).Value = 1508.5  # L31: 1000 -> 1508.5
$ws.Cells.Item(31, 13).Value = -454  # M31: -517 -> -454
$ws.Cells.Item(31, 14).Value = -2004.5  # N31: -1496 -> -2004.5
$ws.Cells.Item(34, 8).Value = 26720.666  # H34: 30574.666 -> 26720.666
$ws.Cells.Item(34, 9).Value = 10325  # I34: 10850 -> 10325
$ws.Cells.Item(34, 10).Value = 59512  # J34: 70024 -> 59512
$ws.Cells.Item(34, 11).Value = 10325  # K34: 10850 -> 10325
$ws.Cells.Item(34, 12).Value = 59512  # L34: 70024 -> 59512
$ws.Cells.Item(34, 13).Value = -10153  # M34: -10678 -> -10153
$ws.Cells.Item(34, 14).Value = -59856  # N34: -70368 -> -59856
$ws.Cells.Item(61, 8).Value = 4406.933  # H61: 4318.9375 -> 4406.933
$ws.Cells.Item(61, 10).Value = 7000.8  # J61: 6333.8335 -> 7000.8
$ws.Cells.Item(61, 12).Value = 7000.8  # L61: 6333.8335 -> 7000.8
$ws.Cells.Item(61, 14).Value = -7404.8  # N61: -6737.8335 -> -7404.8
$ws.Cells.Item(113, 8).Value = 4406.933  # H113: 4318.9375 -> 4406.933
$ws.Cells.Item(113, 10).Value = 7000.8  # J113: 6333.8335 -> 7000.8
$ws.Cells.Item(113, 12).Value = 7000.8  # L113: 6333.8335 -> 7000.8
$ws.Cells.Item(113, 14).Value = -11340.8  # N113: -10673.8335 -> -11340.8
$ws.Cells.Item(122, 8).Value = 1092524  # H122: 936835.9 -> 1092524
$ws.Cells.Item(122, 9).Value = 1403051.9  # I122: 1636737.4 -> 1403051.9
$ws.Cells.Item(122, 10).Value = 5676.25  # J122: 3633.889 -> 5676.25
$ws.Cells.Item(122, 11).Value = 4209155.699999999  # K122: 4910212.199999999 -> 4209155.699999999
$ws.Cells.Item(122, 12).Value = 17028.75  # L122: 10901.667 -> 17028.75
$ws.Cells.Item(122, 13).Value = -4206705.699999999  # M122: -4907762.199999999 -> -4206705.699999999
$ws.Cells.Item(122, 14).Value = -21928.75  # N122: -15801.667 -> -21928.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 5000  # H39: 4808.8 -> 5000
$ws.Cells.Item(39, 9).Value = 0  # I39: 4044 -> 0
$ws.Cells.Item(39, 11).Value = 0  # K39: 4044 -> 0
$ws.Cells.Item(39, 13).ClearContents()  # M39: remove (was -3631)
$ws.Cells.Item(62, 8).Value = 4882.25  # H62: 4852.4 -> 4882.25
$ws.Cells.Item(62, 9).Value = 3277.5  # I62: 3518.3333 -> 3277.5
$ws.Cells.Item(62, 10).Value = 5417.1665  # J62: 5424.143 -> 5417.1665
$ws.Cells.Item(62, 11).Value = 3277.5  # K62: 3518.3333 -> 3277.5
$ws.Cells.Item(62, 12).Value = 5417.1665  # L62: 5424.143 -> 5417.1665
$ws.Cells.Item(62, 13).Value = -2653.5  # M62: -2894.3333 -> -2653.5
$ws.Cells.Item(62, 14).Value = -6665.1665  # N62: -6672.143 -> -6665.1665
$ws.Cells.Item(65, 8).Value = 4882.25  # H65: 4852.4 -> 4882.25
$ws.Cells.Item(65, 9).Value = 3277.5  # I65: 3518.3333 -> 3277.5
$ws.Cells.Item(65, 10).Value = 5417.1665  # J65: 5424.143 -> 5417.1665
$ws.Cells.Item(65, 11).Value = 16387.5  # K65: 17591.6665 -> 16387.5
$ws.Cells.Item(65, 12).Value = 27085.8325  # L65: 27120.715 -> 27085.8325
$ws.Cells.Item(65, 13).Value = -13267.5  # M65: -14471.6665 -> -13267.5
$ws.Cells.Item(65, 14).Value = -33325.8325  # N65: -33360.715 -> -33325.8325
$ws.Cells.Item(122, 8).Value = 1763.9286  # H122: 1867.6 -> 1763.9286
$ws.Cells.Item(122, 9).Value = 1749.6666  # I122: 1908.4117 -> 1749.6666
$ws.Cells.Item(122, 10).Value = 1806.7142  # J122: 1780.875 -> 1806.7142
$ws.Cells.Item(122, 11).Value = 5248.9998  # K122: 5725.2351 -> 5248.9998
$ws.Cells.Item(122, 12).Value = 5420.142599999999  # L122: 5342.625 -> 5420.142599999999
$ws.Cells.Item(122, 13).Value = -2798.9998  # M122: -3275.2351 -> -2798.9998
$ws.Cells.Item(122, 14).Value = -10320.1426  # N122: -10242.625 -> -10320.1426
$ws.Cells.Item(126, 8).Value = 940.3913  # H126: 973.2857 -> 940.3913
$ws.Cells.Item(126, 9).Value = 806.5833  # I126: 843.2222 -> 806.5833
$ws.Cells.Item(126, 10).Value = 1086.3636  # J126: 1070.8334 -> 1086.3636
$ws.Cells.Item(126, 11).Value = 2419.7499  # K126: 2529.6666 -> 2419.7499
$ws.Cells.Item(126, 12).Value = 3259.0908  # L126: 3212.5002 -> 3259.0908
$ws.Cells.Item(126, 13).Value = 50.2501000000002  # M126: -59.66660000000002 -> 50.2501000000002
$ws.Cells.Item(126, 14).Value = -8199.0908  # N126: -8152.5002 -> -8199.0908
$ws.Cells.Item(132, 8).Value = 2388.2  # H132: 1714.7084 -> 2388.2
$ws.Cells.Item(132, 9).Value = 1500  # I132: 1310.4615 -> 1500
$ws.Cells.Item(132, 10).Value = 2610.25  # J132: 2192.4546 -> 2610.25
$ws.Cells.Item(132, 11).Value = 4500  # K132: 3931.3845 -> 4500
$ws.Cells.Item(132, 12).Value = 7830.75  # L132: 6577.3638 -> 7830.75
$ws.Cells.Item(132, 13).Value = -1970  # M132: -1401.3845 -> -1970
$ws.Cells.Item(132, 14).Value = -12890.75  # N132: -11637.3638 -> -12890.75
